$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.212.86'
$ws.Range('E2').Value = '  -0.56%  '
$ws.Range('D3').Value = '3.877.54'
$ws.Range('E3').Value = '  -1.01%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.87'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.82'
$ws.Range('E6').Value = '  +1.44%  '
$ws.Range('D7').Value = '3.877.88'
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.47'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -0.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000250'
$ws.Range('E13').Value = '  +0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.07'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '4.521.91'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '3.874.51'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '68.173.67'
$ws.Range('E17').Value = '  -0.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.18'
$ws.Range('E18').Value = '  +5.70%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.87'
$ws.Range('E21').Value = '  -0.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '470.20'
$ws.Range('E22').Value = '  -3.15%  '
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000164'
$ws.Range('E24').Value = '  -3.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.63'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.17'
$ws.Range('E27').Value = '  +0.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  -0.95%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').Value = '  +0.99%  '
$ws.Range('D31').Value = '4.025.71'
$ws.Range('E31').Value = '  -1.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.95'
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('E33').Value = '  -2.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.39'
$ws.Range('E34').Value = '  -1.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.39'
$ws.Range('E35').Value = '  +1.75%  '
$ws.Range('D36').Value = '3.851.43'
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.105'
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('E38').Value = '  +11.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.02'
$ws.Range('E39').Value = '  -1.70%  '
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.91'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.314'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '436.05'
$ws.Range('E44').Value = '  +1.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.99'
$ws.Range('E45').Value = '  +0.20%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '47.28'
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.58'
$ws.Range('E48').Value = '  +0.94%  '
$ws.Range('E49').Value = '  +7.91%  '
$ws.Range('B50').Value = 'Arweave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '40.56'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '143.68'
$ws.Range('E51').Value = '  +1.32%  '
